$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.9075355
$ws.Range("H2").Value = 33.815071
$ws.Range("I2").Value = 0.2093814122600371
$ws.Range("J2").Value = 0.160920589855456
$ws.Range("M2").Value = 0.071111
$ws.Range("N2").Value = 0.142222
$ws.Range("O2").Value = 0.02711460746047303
$ws.Range("P2").Value = 0.02678527708115022
$ws.Range("Q2").Value = 1.2023117569405
$ws.Range("R2").Value = 4.809247027762
$ws.Range("S2").Value = 0.005677294802950382
$ws.Range("T2").Value = 0.00431030258734052

$ws.Range("G3").Value = 16.9075355
$ws.Range("H3").Value = 33.815071
$ws.Range("I3").Value = 0.2093814122600371
$ws.Range("J3").Value = 0.160920589855456
$ws.Range("O3").Value = 0.0245904030281302
$ws.Range("P3").Value = 0.03643759694506741
$ws.Range("Q3").Value = 1.0903838719305
$ws.Range("R3").Value = 6.542303231583001
$ws.Range("S3").Value = 0.005148773314073394
$ws.Range("T3").Value = 0.005863559593315608

$ws.Range("G4").Value = 16.9075355
$ws.Range("H4").Value = 33.815071
$ws.Range("I4").Value = 0.2093814122600371
$ws.Range("J4").Value = 0.160920589855456
$ws.Range("M4").Value = 2.4870065
$ws.Range("N4").Value = 4.974013
$ws.Range("O4").Value = 0.9482949895113968
$ws.Range("P4").Value = 0.9367771259737823
$ws.Range("Q4").Value = 42.04915068748075
$ws.Range("R4").Value = 168.196602749923
$ws.Range("S4").Value = 0.1985553441430133
$ws.Range("T4").Value = 0.1507467276747999

$ws.Range("I5").Value = 0.005734943705176796
$ws.Range("J5").Value = 0.006611407243820265
$ws.Range("M5").Value = 0.071111
$ws.Range("N5").Value = 0.142222
$ws.Range("O5").Value = 0.02711460746047303
$ws.Range("P5").Value = 0.02678527708115022
$ws.Range("Q5").Value = 0.03293124335966666
$ws.Range("R5").Value = 0.197587460158
$ws.Range("S5").Value = 0.0001555007473737796
$ws.Range("T5").Value = 0.0001770883749220495

$ws.Range("I6").Value = 0.005734943705176796
$ws.Range("J6").Value = 0.006611407243820265
$ws.Range("O6").Value = 0.0245904030281302
$ws.Range("P6").Value = 0.03643759694506741
$ws.Range("S6").Value = 0.0001410245770539357
$ws.Range("T6").Value = 0.0002409037923900218

$ws.Range("I7").Value = 0.005734943705176796
$ws.Range("J7").Value = 0.006611407243820265
$ws.Range("M7").Value = 2.4870065
$ws.Range("N7").Value = 4.974013
$ws.Range("O7").Value = 0.9482949895113968
$ws.Range("P7").Value = 0.9367771259737823
$ws.Range("Q7").Value = 1.151723591126167
$ws.Range("R7").Value = 6.910341546757
$ws.Range("S7").Value = 0.005438418380749081
$ws.Range("T7").Value = 0.006193415076508193

$ws.Range("G8").Value = 7.466036666666668
$ws.Range("H8").Value = 22.39811
$ws.Range("I8").Value = 0.09245873245405202
$ws.Range("J8").Value = 0.106589073045193
$ws.Range("M8").Value = 0.071111
$ws.Range("N8").Value = 0.142222
$ws.Range("O8").Value = 0.02711460746047303
$ws.Range("P8").Value = 0.02678527708115022
$ws.Range("Q8").Value = 0.5309173334033334
$ws.Range("R8").Value = 3.18550400042
$ws.Range("S8").Value = 0.002506982236784519
$ws.Range("T8").Value = 0.002855017855338455

$ws.Range("G9").Value = 7.466036666666668
$ws.Range("H9").Value = 22.39811
$ws.Range("I9").Value = 0.09245873245405202
$ws.Range("J9").Value = 0.106589073045193
$ws.Range("O9").Value = 0.0245904030281302
$ws.Range("P9").Value = 0.03643759694506741
$ws.Range("Q9").Value = 0.4814921706700001
$ws.Range("R9").Value = 4.333429536030001
$ws.Range("S9").Value = 0.002273597494515201
$ws.Range("T9").Value = 0.003883849682369091

$ws.Range("G10").Value = 7.466036666666668
$ws.Range("H10").Value = 22.39811
$ws.Range("I10").Value = 0.09245873245405202
$ws.Range("J10").Value = 0.106589073045193
$ws.Range("M10").Value = 2.4870065
$ws.Range("N10").Value = 4.974013
$ws.Range("O10").Value = 0.9482949895113968
$ws.Range("P10").Value = 0.9367771259737823
$ws.Range("Q10").Value = 18.56808171923834
$ws.Range("R10").Value = 111.40849031543
$ws.Range("S10").Value = 0.0876781527227523
$ws.Range("T10").Value = 0.09985020550748544

$ws.Range("G11").Value = 15.2071285
$ws.Range("H11").Value = 30.414257
$ws.Range("I11").Value = 0.1883237235698756
$ws.Range("J11").Value = 0.1447366523777351
$ws.Range("M11").Value = 0.071111
$ws.Range("N11").Value = 0.142222
$ws.Range("O11").Value = 0.02711460746047303
$ws.Range("P11").Value = 0.02678527708115022
$ws.Range("Q11").Value = 1.0813941147635
$ws.Range("R11").Value = 4.325576459053999
$ws.Range("S11").Value = 0.005106323840091812
$ws.Range("T11").Value = 0.003876811337735754

$ws.Range("G12").Value = 15.2071285
$ws.Range("H12").Value = 30.414257
$ws.Range("I12").Value = 0.1883237235698756
$ws.Range("J12").Value = 0.1447366523777351
$ws.Range("O12").Value = 0.0245904030281302
$ws.Range("P12").Value = 0.03643759694506741
$ws.Range("Q12").Value = 0.9807229240935
$ws.Range("R12").Value = 5.884337544561
$ws.Range("S12").Value = 0.004630956262341425
$ws.Range("T12").Value = 0.005273855802518244

$ws.Range("G13").Value = 15.2071285
$ws.Range("H13").Value = 30.414257
$ws.Range("I13").Value = 0.1883237235698756
$ws.Range("J13").Value = 0.1447366523777351
$ws.Range("M13").Value = 2.4870065
$ws.Range("N13").Value = 4.974013
$ws.Range("O13").Value = 0.9482949895113968
$ws.Range("P13").Value = 0.9367771259737823
$ws.Range("Q13").Value = 37.82022742583525
$ws.Range("R13").Value = 151.280909703341
$ws.Range("S13").Value = 0.1785864434674424
$ws.Range("T13").Value = 0.1355859852374811

$ws.Range("G14").Value = 35.389713
$ws.Range("H14").Value = 106.169139
$ws.Range("I14").Value = 0.4382630506626701
$ws.Range("J14").Value = 0.5052421883817986
$ws.Range("M14").Value = 0.071111
$ws.Range("N14").Value = 0.142222
$ws.Range("O14").Value = 0.02711460746047303
$ws.Range("P14").Value = 0.02678527708115022
$ws.Range("Q14").Value = 2.516597881143
$ws.Range("R14").Value = 15.099587286858
$ws.Range("S14").Value = 0.01188333058314771
$ws.Range("T14").Value = 0.01353305200889317

$ws.Range("G15").Value = 35.389713
$ws.Range("H15").Value = 106.169139
$ws.Range("I15").Value = 0.4382630506626701
$ws.Range("J15").Value = 0.5052421883817986
$ws.Range("O15").Value = 0.0245904030281302
$ws.Range("P15").Value = 0.03643759694506741
$ws.Range("Q15").Value = 2.282317981083
$ws.Range("R15").Value = 20.540861829747
$ws.Range("S15").Value = 0.0107770650481329
$ws.Range("T15").Value = 0.01840981121989979

$ws.Range("G16").Value = 35.389713
$ws.Range("H16").Value = 106.169139
$ws.Range("I16").Value = 0.4382630506626701
$ws.Range("J16").Value = 0.5052421883817986
$ws.Range("M16").Value = 2.4870065
$ws.Range("N16").Value = 4.974013
$ws.Range("O16").Value = 0.9482949895113968
$ws.Range("P16").Value = 0.9367771259737823
$ws.Range("Q16").Value = 88.0144462641345
$ws.Range("R16").Value = 528.086677584807
$ws.Range("S16").Value = 0.4156026550313895
$ws.Range("T16").Value = 0.4732993251530055

$ws.Range("G17").Value = 5.316425333333333
$ws.Range("H17").Value = 15.949276
$ws.Range("I17").Value = 0.06583813734818843
$ws.Range("J17").Value = 0.07590008909599709
$ws.Range("M17").Value = 0.071111
$ws.Range("N17").Value = 0.142222
$ws.Range("O17").Value = 0.02711460746047303
$ws.Range("P17").Value = 0.02678527708115022
$ws.Range("Q17").Value = 0.3780563218786666
$ws.Range("R17").Value = 2.268337931272
$ws.Range("S17").Value = 0.001785175250124838
$ws.Range("T17").Value = 0.002033004916920271

$ws.Range("G18").Value = 5.316425333333333
$ws.Range("H18").Value = 15.949276
$ws.Range("I18").Value = 0.06583813734818843
$ws.Range("J18").Value = 0.07590008909599709
$ws.Range("O18").Value = 0.0245904030281302
$ws.Range("P18").Value = 0.03643759694506741
$ws.Range("Q18").Value = 0.342861586172
$ws.Range("R18").Value = 3.085754275548
$ws.Range("S18").Value = 0.001618986332013345
$ws.Range("T18").Value = 0.002765616854574648

$ws.Range("G19").Value = 5.316425333333333
$ws.Range("H19").Value = 15.949276
$ws.Range("I19").Value = 0.06583813734818843
$ws.Range("J19").Value = 0.07590008909599709
$ws.Range("M19").Value = 2.4870065
$ws.Range("N19").Value = 4.974013
$ws.Range("O19").Value = 0.9482949895113968
$ws.Range("P19").Value = 0.9367771259737823
$ws.Range("Q19").Value = 13.22198436076467
$ws.Range("R19").Value = 79.331906164588
$ws.Range("S19").Value = 0.06243397576605025
$ws.Range("T19").Value = 0.07110146732450216
